$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 31622
$ws.Range("E2").Value = 1495
$ws.Range("F2").Value = 1495
$ws.Range("G2").Value = 1671
$ws.Range("H2").Value = 1228
$ws.Range("I2").Value = 736
$ws.Range("J2").Value = 492
$ws.Range("K2").Value = 38796
$ws.Range("L2").Value = 10265
$ws.Range("M2").Value = 28531
$ws.Range("N2").Value = 24823
$ws.Range("O2").Value = 3708
$ws.Range("P2").Value = 56
$ws.Range("Q2").Value = 3119
$ws.Range("R2").Value = -3020
$ws.Range("S2").Value = -1736
$ws.Range("T2").Value = 1595
$ws.Range("U2").Value = 1524
$ws.Range("V2").Value = 3378
$ws.Range("W2").Value = 4.73
$ws.Range("X2").Value = 3.88
$ws.Range("Y2").Value = 2.98
$ws.Range("Z2").Value = 3.12
$ws.Range("AA2").Value = 35.98
$ws.Range("AB2").Value = 47279.15
$ws.Range("AC2").Value = 66136
$ws.Range("AD2").Value = 17.48
$ws.Range("AE2").Value = 2949412
$ws.Range("AF2").Value = 0.39
$ws.Range("AG2").Value = 1250
$ws.Range("AH2").Value = 0.11
$ws.Range("AI2").Value = 1.43
$ws.Range("AJ2").Value = 1113400

# Row 3
$ws.Range("D3").Value = 28044
$ws.Range("E3").Value = 1598
$ws.Range("F3").Value = 1598
$ws.Range("G3").Value = 1555
$ws.Range("H3").Value = 1046
$ws.Range("I3").Value = 750
$ws.Range("J3").Value = 296
$ws.Range("K3").Value = 38249
$ws.Range("L3").Value = 8748
$ws.Range("M3").Value = 29501
$ws.Range("N3").Value = 25396
$ws.Range("O3").Value = 4105
$ws.Range("P3").Value = 56
$ws.Range("Q3").Value = 3734
$ws.Range("R3").Value = -2216
$ws.Range("S3").Value = -1410
$ws.Range("T3").Value = 1503
$ws.Range("U3").Value = 2231
$ws.Range("V3").Value = 2132
$ws.Range("W3").Value = 5.7
$ws.Range("X3").Value = 3.73
$ws.Range("Y3").Value = 2.99
$ws.Range("Z3").Value = 2.72
$ws.Range("AA3").Value = 29.65
$ws.Range("AB3").Value = 48134.66
$ws.Range("AC3").Value = 67395
$ws.Range("AD3").Value = 16.17
$ws.Range("AE3").Value = 3017474
$ws.Range("AF3").Value = 0.36
$ws.Range("AG3").Value = 1750
$ws.Range("AH3").Value = 0.16
$ws.Range("AI3").Value = 1.96
$ws.Range("AJ3").Value = 1113400

# Row 4
$ws.Range("D4").Value = 26711
$ws.Range("E4").Value = 1601
$ws.Range("F4").Value = 1601
$ws.Range("G4").Value = 849
$ws.Range("H4").Value = 445
$ws.Range("I4").Value = 143
$ws.Range("J4").Value = 302
$ws.Range("K4").Value = 38223
$ws.Range("L4").Value = 8451
$ws.Range("M4").Value = 29772
$ws.Range("N4").Value = 25427
$ws.Range("O4").Value = 4346
$ws.Range("P4").Value = 56
$ws.Range("Q4").Value = 3487
$ws.Range("R4").Value = -3091
$ws.Range("S4").Value = -999
$ws.Range("T4").Value = 828
$ws.Range("U4").Value = 2659
$ws.Range("V4").Value = 1213
$ws.Range("W4").Value = 5.99
$ws.Range("X4").Value = 1.67
$ws.Range("Y4").Value = 0.5600000000000001
$ws.Range("Z4").Value = 1.17
$ws.Range("AA4").Value = 28.38
$ws.Range("AB4").Value = 48351.59
$ws.Range("AC4").Value = 12841
$ws.Range("AD4").Value = 73.75
$ws.Range("AE4").Value = 3021101
$ws.Range("AF4").Value = 0.31
$ws.Range("AG4").Value = 1750
$ws.Range("AH4").Value = 0.18
$ws.Range("AI4").Value = 10.3
$ws.Range("AJ4").Value = 1113400

# Row 5
$ws.Range("D5").Value = 29158
$ws.Range("E5").Value = 2412
$ws.Range("F5").Value = 2412
$ws.Range("G5").Value = 2758
$ws.Range("H5").Value = 1792
$ws.Range("I5").Value = 1388
$ws.Range("J5").Value = 404
$ws.Range("K5").Value = 41095
$ws.Range("L5").Value = 9352
$ws.Range("M5").Value = 31743
$ws.Range("N5").Value = 27063
$ws.Range("O5").Value = 4680
$ws.Range("P5").Value = 56
$ws.Range("Q5").Value = 3467
$ws.Range("R5").Value = -3054
$ws.Range("S5").Value = -93
$ws.Range("T5").Value = 685
$ws.Range("U5").Value = 2782
$ws.Range("V5").Value = 1005
$ws.Range("W5").Value = 8.27
$ws.Range("X5").Value = 6.14
$ws.Range("Y5").Value = 5.29
$ws.Range("Z5").Value = 4.52
$ws.Range("AA5").Value = 29.46
$ws.Range("AB5").Value = 50769.36
$ws.Range("AC5").Value = 124662
$ws.Range("AD5").Value = 10.56
$ws.Range("AE5").Value = 3215563
$ws.Range("AF5").Value = 0.41
$ws.Range("AG5").Value = 1925
$ws.Range("AH5").Value = 0.15
$ws.Range("AI5").Value = 1.17
$ws.Range("AJ5").Value = 1113400

# Row 6
$ws.Range("D6").Value = 31088
$ws.Range("E6").Value = 3315
$ws.Range("F6").Value = 3315
$ws.Range("G6").Value = 3438
$ws.Range("H6").Value = 2491
$ws.Range("I6").Value = 2141
$ws.Range("K6").Value = 45027
$ws.Range("L6").Value = 9977
$ws.Range("M6").Value = 35050
$ws.Range("N6").Value = 30863
$ws.Range("P6").Value = 56
$ws.Range("Q6").Value = 3581
$ws.Range("R6").Value = -3870
$ws.Range("S6").Value = -163
$ws.Range("T6").Value = 817
$ws.Range("U6").Value = 2764
$ws.Range("V6").Value = 970
$ws.Range("W6").Value = 10.66
$ws.Range("X6").Value = 8.01
$ws.Range("Y6").Value = 7.39
$ws.Range("Z6").Value = 5.79
$ws.Range("AA6").Value = 28.46
$ws.Range("AB6").Value = 56497.53
$ws.Range("AC6").Value = 192300
$ws.Range("AD6").Value = 6.86
$ws.Range("AE6").Value = 3667047
$ws.Range("AF6").Value = 0.36
$ws.Range("AG6").Value = 3000
$ws.Range("AH6").Value = 0.23
$ws.Range("AI6").Value = 1.18
$ws.Range("AJ6").Value = 1113400

# Row 7 - clear remaining data cells
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8 - clear remaining data cells
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9 - clear remaining data cells
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
